$wb = $excel.ActiveWorkbook

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" everywhere it appears ---
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# --- Widen the "Status"/"zh-cn"/"de-de" columns to fit the longer text ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527

# --- Widen "Latest Target File" / "Latest Handback File" columns (I, J) ---
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40

# --- zh-cn sheet: populate handback report columns ---
$wsZh.Range("I2").Value = "4c94b8e1-af2a-4aa5-a24b-45c379346fef.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/276a272a9e959a172a85c7d12edc4ebc1ebbf8b1/e2e/4c94b8e1-af2a-4aa5-a24b-45c379346fef.md", "", "", "4c94b8e1-af2a-4aa5-a24b-45c379346fef.md") | Out-Null
$wsZh.Range("J2").Value = "4c94b8e1-af2a-4aa5-a24b-45c379346fef.474a4e8d368245d15003d7553a9e3855ea771211.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-19 17:05:30"

$wsZh.Range("I3").Value = "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/276a272a9e959a172a85c7d12edc4ebc1ebbf8b1/e2e/b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md", "", "", "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md") | Out-Null
$wsZh.Range("J3").Value = "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.fa45dbba1a061242178c26dc46c2609a48b9bb04.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-19 17:05:30"

# --- de-de sheet: populate handback report columns ---
$wsDe.Range("I2").Value = "4c94b8e1-af2a-4aa5-a24b-45c379346fef.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/276a272a9e959a172a85c7d12edc4ebc1ebbf8b1/e2e/4c94b8e1-af2a-4aa5-a24b-45c379346fef.md", "", "", "4c94b8e1-af2a-4aa5-a24b-45c379346fef.md") | Out-Null
$wsDe.Range("J2").Value = "4c94b8e1-af2a-4aa5-a24b-45c379346fef.474a4e8d368245d15003d7553a9e3855ea771211.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-19 17:05:38"

$wsDe.Range("I3").Value = "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/276a272a9e959a172a85c7d12edc4ebc1ebbf8b1/e2e/b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md", "", "", "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.md") | Out-Null
$wsDe.Range("J3").Value = "b630bac9-0c6d-49b6-93dd-05b39a6f20ca.fa45dbba1a061242178c26dc46c2609a48b9bb04.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-19 17:05:38"
